# Update "paises.xlsx" country/COVID stats and refresh timestamp.
# Country order (column A) is driven by the shared-string table upstream;
# here we just overwrite each row's country name + stats to match the new
# snapshot. Several countries with close case counts swapped adjacent rows
# (Argentina/Banglades, Zambia/Guayana Francesa, Surinam/Mali/Mozambique,
# Vietnam/Santo Tome y Principe, Islas Malvinas/Montserrat).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 02:10"

$ws.Range("B4").Value = 5359563
$ws.Range("C4").Value = 53606
$ws.Range("D4").Value = 2804232
$ws.Range("E4").Value = 2386219
$ws.Range("G4").Value = 1367
$ws.Range("H4").Value = 169112

$ws.Range("B5").Value = 3170474
$ws.Range("C5").Value = 58081
$ws.Range("E5").Value = 756734
$ws.Range("G5").Value = 1164
$ws.Range("H5").Value = 104263

$ws.Range("A18").Value = "Argentina"
$ws.Range("B18").Value = 268574
$ws.Range("C18").Value = 7663
$ws.Range("D18").Value = 187283
$ws.Range("E18").Value = 76078
$ws.Range("G18").Value = 209
$ws.Range("H18").Value = 5213

$ws.Range("A19").Value = "Banglades"
$ws.Range("B19").Value = 266498
$ws.Range("C19").Value = 2995
$ws.Range("D19").Value = 153089
$ws.Range("E19").Value = 109896
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 3513

$ws.Range("B27").Value = 120844
$ws.Range("C27").Value = 423
$ws.Range("D27").Value = 107148
$ws.Range("E27").Value = 4690
$ws.Range("G27").Value = 15
$ws.Range("H27").Value = 9006

$ws.Range("B46").Value = 59089
$ws.Range("C46").Value = 1123
$ws.Range("D46").Value = 47394
$ws.Range("E46").Value = 9428
$ws.Range("G46").Value = 34
$ws.Range("H46").Value = 2267

$ws.Range("B50").Value = 50210
$ws.Range("C50").Value = 1282
$ws.Range("D50").Value = 34888
$ws.Range("E50").Value = 14263
$ws.Range("G50").Value = 7
$ws.Range("H50").Value = 1059

$ws.Range("B74").Value = 19075
$ws.Range("C74").Value = 292
$ws.Range("E74").Value = 5277

$ws.Range("B75").Value = 18263
$ws.Range("C75").Value = 50
$ws.Range("E75").Value = 2542
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 401

$ws.Range("A89").Value = "Zambia"
$ws.Range("B89").Value = 8501
$ws.Range("C89").Value = 226
$ws.Range("D89").Value = 7233
$ws.Range("E89").Value = 1022
$ws.Range("G89").Value = 5
$ws.Range("H89").Value = 246

$ws.Range("A90").Value = "Guayana Francesa"
$ws.Range("B90").Value = 8423
$ws.Range("C90").Value = 63
$ws.Range("D90").Value = 7713
$ws.Range("E90").Value = 660
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 50

$ws.Range("B94").Value = 7743
$ws.Range("C94").Value = 94
$ws.Range("D94").Value = 5123
$ws.Range("E94").Value = 2433
$ws.Range("G94").Value = 4
$ws.Range("H94").Value = 187

$ws.Range("B98").Value = 7300
$ws.Range("C98").Value = 58
$ws.Range("D98").Value = 6262
$ws.Range("E98").Value = 916

$ws.Range("B100").Value = 6622
$ws.Range("C100").Value = 24
$ws.Range("D100").Value = 5741
$ws.Range("E100").Value = 724

$ws.Range("B113").Value = 3813
$ws.Range("C113").Value = 65
$ws.Range("D113").Value = 2638
$ws.Range("E113").Value = 1102
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 73

$ws.Range("B119").Value = 3128
$ws.Range("C119").Value = 35
$ws.Range("D119").Value = 2504
$ws.Range("E119").Value = 536

$ws.Range("A124").Value = "Surinam"
$ws.Range("B124").Value = 2653
$ws.Range("C124").Value = 94
$ws.Range("D124").Value = 1789
$ws.Range("E124").Value = 825
$ws.Range("H124").Value = 39

$ws.Range("A125").Value = "Mali"
$ws.Range("B125").Value = 2582
$ws.Range("C125").Value = 5
$ws.Range("D125").Value = 1977
$ws.Range("E125").Value = 480
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 125

$ws.Range("A126").Value = "Mozambique"
$ws.Range("C126").Value = 78
$ws.Range("D126").Value = 951
$ws.Range("E126").Value = 1589
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 19

$ws.Range("B151").Value = 1161
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 1075
$ws.Range("E151").Value = 17

$ws.Range("B155").Value = 1036
$ws.Range("C155").Value = 47
$ws.Range("D155").Value = 122
$ws.Range("E155").Value = 899

$ws.Range("A158").Value = "Vietnam"
$ws.Range("B158").Value = 883
$ws.Range("C158").Value = 17
$ws.Range("D158").Value = 409
$ws.Range("E158").Value = 457
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 17

$ws.Range("A159").Value = "Santo Tome y Principe"
$ws.Range("B159").Value = 882
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 807
$ws.Range("E159").Value = 60
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 15

$ws.Range("B161").Value = 798
$ws.Range("C161").Value = 17
$ws.Range("E161").Value = 599

$ws.Range("B165").Value = 623
$ws.Range("C165").Value = 21
$ws.Range("D165").Value = 191
$ws.Range("E165").Value = 410

$ws.Range("B168").Value = 409
$ws.Range("C168").Value = 1
$ws.Range("E168").Value = 93

$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
